# Update the two-digit ÷ one-digit division problems to new values.
# Each original expression is unique in the document, so a simple
# Find/Replace (MatchCase, no wildcards) targets the correct cell.

$d = $word.ActiveDocument

$replacements = @(
    @("75÷3=", "79÷8="),
    @("25÷5=", "45÷7="),
    @("43÷3=", "48÷4="),
    @("10÷2=", "14÷5="),
    @("33÷6=", "89÷4="),
    @("49÷5=", "58÷5="),
    @("52÷7=", "24÷2="),
    @("11÷8=", "12÷2="),
    @("80÷6=", "60÷9="),
    @("20÷6=", "13÷4="),
    @("71÷6=", "50÷3="),
    @("30÷8=", "84÷9="),
    @("35÷3=", "49÷5="),
    @("15÷4=", "62÷5="),
    @("21÷8=", "21÷2="),
    @("78÷2=", "26÷2="),
    @("69÷5=", "72÷7="),
    @("87÷6=", "67÷8="),
    @("71÷5=", "97÷3="),
    @("35÷9=", "87÷2="),
    @("62÷2=", "14÷3="),
    @("89÷2=", "18÷5="),
    @("72÷5=", "11÷5="),
    @("67÷9=", "74÷2="),
    @("25÷9=", "61÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
